$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.885.42'
$ws.Range('E2').Value = '  -1.21%  '
$ws.Range('D3').Value = '2.330.68'
$ws.Range('E3').Value = '  +1.05%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.23'
$ws.Range('E5').Value = '  -1.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.37'
$ws.Range('E6').Value = '  -2.69%  '
$ws.Range('E7').Value = '  -4.21%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.506'
$ws.Range('E9').Value = '  -4.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.49'
$ws.Range('E10').Value = '  -5.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.28'
$ws.Range('E11').Value = '  +1.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0793'
$ws.Range('E12').Value = '  -2.23%  '
$ws.Range('E13').Value = '  +0.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.76'
$ws.Range('E14').Value = '  -4.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.78'
$ws.Range('E15').Value = '  +4.66%  '
$ws.Range('D16').Value = '2.361.80'
$ws.Range('E16').Value = '  +2.64%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.826'
$ws.Range('E17').Value = '  +1.90%  '
$ws.Range('D18').Value = '42.833.67'
$ws.Range('E18').Value = '  -1.09%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.15'
$ws.Range('E19').Value = '  -0.55%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0904'
$ws.Range('E20').Value = '  -2.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.58'
$ws.Range('E21').Value = '  -5.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.19'
$ws.Range('E22').Value = '  +1.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.65'
$ws.Range('E23').Value = '  -2.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.98'
$ws.Range('E24').Value = '  -2.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.36'
$ws.Range('E27').Value = '  +2.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.95'
$ws.Range('E28').Value = '  -0.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.18'
$ws.Range('E29').Value = '  -5.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.68'
$ws.Range('E30').Value = '  -6.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.22'
$ws.Range('E31').Value = '  -4.62%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '160.83'
$ws.Range('E32').Value = '  -4.59%  '
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.06'
$ws.Range('E34').Value = '  -4.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.59'
$ws.Range('E35').Value = '  +2.69%  '
$ws.Range('E36').Value = '  -3.49%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0718'
$ws.Range('E37').Value = '  -3.38%  '
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '17.12'
$ws.Range('E38').Value = '  -5.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.83'
$ws.Range('E39').Value = '  -2.48%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.88'
$ws.Range('E40').Value = '  -5.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.101'
$ws.Range('E41').Value = '  -4.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.112'
$ws.Range('E42').Value = '  -3.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.57'
$ws.Range('E43').Value = '  -0.07%  '
$ws.Range('D44').Value = '2.011.61'
$ws.Range('E44').Value = '  +1.27%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0282'
$ws.Range('E45').Value = '  -4.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '18.64'
$ws.Range('E46').Value = '  -1.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.17'
$ws.Range('E47').Value = '  +1.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.88'
$ws.Range('E48').Value = '  -4.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '55.41'
$ws.Range('E49').Value = '  -1.14%  '
$ws.Range('E50').Value = '  -2.11%  '
$ws.Range('D51').Value = '2.561.64'
$ws.Range('E51').Value = '  +1.14%  '

Write-Host "Applied 97 cell updates"
